$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.145.70"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.36"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.88"
$ws.Range("E5").Value = "  +2.96%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.84"
$ws.Range("E8").Value = "  +5.06%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0989"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.121.34"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.38"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.839.06"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.143.93"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.92"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.61"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.18"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.08"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.82"
$ws.Range("E27").Value = "  +19.56%  "
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +28.63%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.00"
$ws.Range("E35").Value = "  +9.48%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.808"
$ws.Range("E36").Value = "  +15.55%  "
$ws.Range("E37").Value = "  +4.55%  "
$ws.Range("E38").Value = "  +7.96%  "
$ws.Range("E39").Value = "  +3.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "89.75"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.340.44"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.97"
$ws.Range("E43").Value = "  +53.20%  "
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.74"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0553"
$ws.Range("E47").Value = "  +6.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.45"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.039.82"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0676"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").Value = "  +0.44%  "
